$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose D/M/N/O/P/R/S values get reshuffled among each other.
# (Row 11 is excluded - it keeps its original values.)
$rows = @(2,3,4,5,6,7,8,9,10,12,13,14,15,16,17,18,19)

# Snapshot the "before" values for every affected row, per column,
# since the permutation below reassigns values between rows and some
# of the writes would otherwise clobber a value still needed later.
$D = @{}
$M = @{}
$N = @{}
$O = @{}
$P = @{}
$R = @{}
$S = @{}

foreach ($r in $rows) {
    $D[$r] = $ws.Cells.Item($r, 4).Value2
    $M[$r] = $ws.Cells.Item($r, 13).Value2
    $N[$r] = $ws.Cells.Item($r, 14).Value2
    $O[$r] = $ws.Cells.Item($r, 15).Value2
    $P[$r] = $ws.Cells.Item($r, 16).Value2
    $R[$r] = $ws.Cells.Item($r, 18).Value2
    $S[$r] = $ws.Cells.Item($r, 19).Value2
}

# Destination row -> source row (which row's original D/M/N/O/P/R/S
# values now land on the destination row).
$map = @{
    2  = 18
    3  = 9
    4  = 16
    5  = 4
    6  = 5
    7  = 14
    8  = 3
    9  = 10
    10 = 8
    12 = 6
    13 = 15
    14 = 17
    15 = 12
    16 = 7
    17 = 19
    18 = 13
    19 = 2
}

foreach ($dst in $rows) {
    $src = $map[$dst]
    $ws.Cells.Item($dst, 4).Value  = $D[$src]
    $ws.Cells.Item($dst, 13).Value = $M[$src]
    $ws.Cells.Item($dst, 14).Value = $N[$src]
    $ws.Cells.Item($dst, 15).Value = $O[$src]
    $ws.Cells.Item($dst, 16).Value = $P[$src]
    $ws.Cells.Item($dst, 18).Value = $R[$src]
    $ws.Cells.Item($dst, 19).Value = $S[$src]
}
